# Add 2022-Q1 data:
#  - insert a new "2022-Q1" worksheet (holdings detail) right before the
#    "总计" (totals) sheet
#  - prepend a 2022-Q1 summary row to the "总计" sheet

function Set-TextCell($cell, $val) {
    # Force the cell to be stored as text even when the literal looks
    # numeric (e.g. "007835", "9.37"), without leaving a lingering
    # text number-format behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$originalActive = $wb.ActiveSheet

$totalSheetBefore = $wb.Worksheets.Item($wb.Worksheets.Count)
$templateSheet = $wb.Worksheets.Item($wb.Worksheets.Count - 1)

# --- 1. New "2022-Q1" holdings sheet, inserted right before "总计" ---
# NOTE: Worksheets.Add(before) inserts the new sheet into the slot that
# $before used to occupy and bumps the old sheet after it — so the
# $totalSheetBefore handle no longer tracks the "总计" sheet once Add()
# returns. Re-resolve it by name afterwards for every later write.
$q1 = $wb.Worksheets.Add($totalSheetBefore)
$q1.Name = "2022-Q1"

$q1.Cells.Item(1, 2).Value = "基金代码"
$q1.Cells.Item(1, 3).Value = "基金名称"
$q1.Cells.Item(1, 4).Value = "基金规模"
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$q1.Cells.Item(1, 6).Value = "仓位占比"
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1.Cells.Item(1, 8).Value = "仓位排名"

# Match the header + index-column styling used on the other quarterly
# sheets (bold, centered, thin border).
$templateSheet.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$templateSheet.Range("A2").Copy()
$q1.Range("A2").PasteSpecial(-4122)

$q1.Cells.Item(2, 1).Value = 0
Set-TextCell $q1.Cells.Item(2, 2) "007835"
Set-TextCell $q1.Cells.Item(2, 3) "国泰鑫睿混合"
Set-TextCell $q1.Cells.Item(2, 4) "9.37"
Set-TextCell $q1.Cells.Item(2, 5) "78.94"
Set-TextCell $q1.Cells.Item(2, 6) "3.05"
Set-TextCell $q1.Cells.Item(2, 7) "0.2858"
$q1.Cells.Item(2, 8).Value = 9

$q1.PageSetup.LeftMargin = $templateSheet.PageSetup.LeftMargin
$q1.PageSetup.RightMargin = $templateSheet.PageSetup.RightMargin
$q1.PageSetup.TopMargin = $templateSheet.PageSetup.TopMargin
$q1.PageSetup.BottomMargin = $templateSheet.PageSetup.BottomMargin
$q1.PageSetup.HeaderMargin = $templateSheet.PageSetup.HeaderMargin
$q1.PageSetup.FooterMargin = $templateSheet.PageSetup.FooterMargin

# --- 2. Prepend the 2022-Q1 row to the "总计" sheet, pushing the rest down ---
$tot = $wb.Worksheets.Item("总计")

# Row 7 is brand new (the sheet used to stop at row 6) so its index-column
# cell needs the bordered "A column" style copied over explicitly; rows
# 2-6 already carry that style on column A and keep it across the Value
# re-writes below.
$tot.Range("A6").Copy()
$tot.Range("A7").PasteSpecial(-4122)

$tot.Cells.Item(7, 1).Value = 5
$tot.Cells.Item(7, 2).Value = "2020-Q4"
$tot.Cells.Item(7, 3).Value = 1
$tot.Cells.Item(7, 4).Value = 1.46

$tot.Cells.Item(6, 1).Value = 4
$tot.Cells.Item(6, 2).Value = "2021-Q1"
$tot.Cells.Item(6, 3).Value = 3
$tot.Cells.Item(6, 4).Value = 1.64

$tot.Cells.Item(5, 1).Value = 3
$tot.Cells.Item(5, 2).Value = "2021-Q2"
$tot.Cells.Item(5, 3).Value = 6
$tot.Cells.Item(5, 4).Value = 5.7

$tot.Cells.Item(4, 1).Value = 2
$tot.Cells.Item(4, 2).Value = "2021-Q3"
$tot.Cells.Item(4, 3).Value = 2
$tot.Cells.Item(4, 4).Value = 0.29

$tot.Cells.Item(3, 1).Value = 1
$tot.Cells.Item(3, 2).Value = "2021-Q4"
$tot.Cells.Item(3, 3).Value = 1
$tot.Cells.Item(3, 4).Value = 0.35

$tot.Cells.Item(2, 1).Value = 0
$tot.Cells.Item(2, 2).Value = "2022-Q1"
$tot.Cells.Item(2, 3).Value = 1
$tot.Cells.Item(2, 4).Value = 0.29

# Restore the tab selection to whatever was active before we added sheets.
$originalActive.Activate()
